# "created page 6 and quit button"
# - Update Saldo (E2:E4) on the "database user" sheet
# - Append 6 new log rows (33-38) to the "log" sheet

$wb = $excel.ActiveWorkbook

# --- database user: update Saldo values -------------------------------
$wsUser = $wb.Worksheets.Item("database user")
$wsUser.Range("E2").Value = 14610
$wsUser.Range("E3").Value = 48760
$wsUser.Range("E4").Value = 199290

# --- log: append new trip rows -----------------------------------------
$wsLog = $wb.Worksheets.Item("log")

$rows = @(
    @("20:38:24", "Taman Mini",    4.5,  "20:38:24", "Tambak Sumur",   5,    806.5, "bagas@mail.com", 1, "Jakarta",  "Surabaya", 693850),
    @("20:41:52", "Semarang",      0,    "20:41:52", "Taman Mini",     4.5,  450.5, "bagas@mail.com", 1, "Semarang", "Jakarta",  358150),
    @("20:47:04", "Juanda",        12.8, "20:47:04", "Serpong",        10.1, 819.9, "bagas@mail.com", 2, "Surabaya", "Jakarta",  1052940),
    @("20:55:04", "Semarang",      0,    "20:55:04", "Bawen",          23.1, 23.1,  "test1@mail.com", 2, "Semarang", "Semarang", 9240),
    @("21:06:11", "Juanda",        12.8, "21:06:11", "Tanjung Priok",  12.1, 821.9, "zaki@mail.com",  1, "Surabaya", "Jakarta",  699640),
    @("21:09:58", "Bawen",         23.1, "21:09:58", "Solo",           40,   63.1,  "zaki@mail.com",  2, "Semarang", "Semarang", 25240)
)

$startRow = 33
for ($i = 0; $i -lt $rows.Count; $i++) {
    $r = $startRow + $i
    $data = $rows[$i]
    for ($c = 0; $c -lt $data.Count; $c++) {
        $wsLog.Cells.Item($r, $c + 1).Value = $data[$c]
    }
}
